$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, centered, bordered) from H1 onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-31 for columns I (I0) and J (IF)
$data = @{
    2  = @(6, 7)
    3  = @(6, 6)
    4  = @(5, 6)
    5  = @(11, 11)
    6  = @(8, 9)
    7  = @(1, 2)
    8  = @(1, 4)
    9  = @(6, 8)
    10 = @(1, 4)
    11 = @(1, 3)
    12 = @(1, 6)
    13 = @(1, 6)
    14 = @(1, 6)
    15 = @(1, 7)
    16 = @(1, 5)
    17 = @(1, 5)
    18 = @(1, 4)
    19 = @(1, 6)
    20 = @(1, 3)
    21 = @(1, 5)
    22 = @(1, 6)
    23 = @(1, 4)
    24 = @(1, 6)
    25 = @(1, 4)
    26 = @(1, 5)
    27 = @(1, 6)
    28 = @(1, 3)
    29 = @(4, 5)
    30 = @(1, 3)
    31 = @(1, 2)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
